$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Daily attendance processing refresh: reorder "Recorded By" email lists
# (order-insensitive sets, re-emitted by the generator in a new order),
# update counts/coverage now that session "HISTOLOGY C1 #2" (row 10) has
# been recorded, and refresh the dependent summary figures.
# ---------------------------------------------------------------------------

# --- G2 / G3: recipient list re-ordered (same people, new order) ----------
$ws.Range("G2").Value = "System, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# --- Class Statistics block (K4:L10) ---------------------------------------
# Recorded Sessions 7 -> 8, Missing Sessions 2 -> 1
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 1

# Coverage % text cell (stored as literal text, not a numeric percentage)
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "27.6%"
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

# --- Row 10: HISTOLOGY C1 session #2 flips from Not Recorded -> Recorded --
# Re-style the whole row to the "Recorded" look (copy format from row 2,
# which already carries that exact style) then fix up the values that
# changed.
$ws.Range("A2:I2").Copy()
$ws.Range("A10:I10").PasteSpecial(-4122)

$ws.Range("G10").Value = "Safa.hany@med.asu.edu.eg"
$ws.Range("H10").Value = "8/251"
$ws.Range("I10").Value = "Recorded"

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "25.0%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# --- G15: recipient list re-ordered ----------------------------------------
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# --- Row 15: PARASITOLOGY group-statistics summary row ---------------------
$ws.Range("O15").Value = 8
$ws.Range("P15").Value = 1

$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "27.6%"
$ws.Range("Q15").Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "25.0%"
$ws.Range("R15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

# --- G28: recipient list re-ordered -----------------------------------------
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
